$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2 and 3 (Tretåig hackspett / Trådticka records) were
# reordered: row 2 now holds the record that used to be on row 3, and
# row 3 now holds the record that used to be on row 2. Swap every cell
# that actually differs between the two rows (columns A, B, D-H and the
# sparse "Aktivitet"/M column), leaving the remaining, identical columns
# untouched.

$columns = @("A", "B", "D", "E", "F", "G", "H", "M")

foreach ($col in $columns) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $val2 = $cell2.Value2
    $val3 = $cell3.Value2

    $cell2.Value2 = $val3
    $cell3.Value2 = $val2
}
